$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("EVCRSbRIC")

# Insert new columns from right to left so earlier column letters remain valid.

# "ISIC 35T39" (col X) needs to become 3 columns: ISIC 351, ISIC 352T353, ISIC 36T39
$ws.Columns("Y:Z").Insert()

# "ISIC 24" (col P) needs to become 2 columns: ISIC 241, ISIC 242
$ws.Columns("Q:Q").Insert()

# "ISIC 23" (col O) needs to become 2 columns: ISIC 231, ISIC 239
$ws.Columns("P:P").Insert()

# Relabel headers (row 1) for the split columns
$ws.Range("O1").Value = "ISIC 231"
$ws.Range("P1").Value = "ISIC 239"
$ws.Range("Q1").Value = "ISIC 241"
$ws.Range("R1").Value = "ISIC 242"
$ws.Range("Z1").Value = "ISIC 351"
$ws.Range("AA1").Value = "ISIC 352T353"
$ws.Range("AB1").Value = "ISIC 36T39"

# New columns in row 2 should carry the same "0" value as their siblings
$ws.Range("P2").Value = 0
$ws.Range("Q2").Value = 0
$ws.Range("R2").Value = 0
$ws.Range("AA2").Value = 0
$ws.Range("AB2").Value = 0
